$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O2").Value = 2937
$ws.Range("O3").Value = 2937
$ws.Range("O4").Value = 1910
$ws.Range("O5").Value = 2937
$ws.Range("O8").Value = 3605
$ws.Range("O9").Value = 4646
$ws.Range("O10").Value = 3605
$ws.Range("O11").Value = 1009
